# Update the "Relative Species Richness" row (row 3) to the new
# "Rarefied Richness" naming used for the var_alias_lookup table.
#
# Old:
#   B3 = rich_tottree
#   C3 = Relative Species Richness
#   D3 = Richness normalized by total trees
# New:
#   B3 = rare_rich
#   C3 = Rarefied Richness
#   D3 = Rarefied Species Richness

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Rarefied Richness"
$ws.Range("D3").Value = "Rarefied Species Richness"
$ws.Range("B3").Value = "rare_rich"

# Move the active selection to B3, matching the saved workbook state.
$ws.Range("B3").Select()
